$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 19500
$ws.Range("I54").Value = 19500
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 19500
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -19014
$ws.Range("N54").ClearContents()
$ws.Range("H76").Value = 4144.095
$ws.Range("I76").Value = 3940.75
$ws.Range("K76").Value = 3940.75
$ws.Range("M76").Value = -3625.75
$ws.Range("H79").Value = 4144.095
$ws.Range("I79").Value = 3940.75
$ws.Range("K79").Value = 3940.75
$ws.Range("M79").Value = -2848.75
$ws.Range("H118").Value = 749.3077
$ws.Range("I118").Value = 703.4167
$ws.Range("K118").Value = 2110.2501
$ws.Range("M118").Value = -453.2501000000002
$ws.Range("H137").Value = 8996.75
$ws.Range("I137").Value = 11612.419
$ws.Range("J137").Value = 2759.3845
$ws.Range("K137").Value = 34837.257
$ws.Range("L137").Value = 8278.1535
$ws.Range("M137").Value = -32287.257
$ws.Range("N137").Value = -13378.1535
$ws.Range("H138").Value = 15767.527
$ws.Range("I138").Value = 1381.7637
$ws.Range("K138").Value = 4145.2911
$ws.Range("M138").Value = 994.7088999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28953.078
$ws.Range("I32").Value = 35834.066
$ws.Range("K32").Value = 35834.066
$ws.Range("M32").Value = -35547.066
$ws.Range("H45").Value = 2559.2
$ws.Range("I45").Value = 1834.4117
$ws.Range("K45").Value = 1834.4117
$ws.Range("M45").Value = -1457.4117
$ws.Range("H61").Value = 6303.136
$ws.Range("I61").Value = 1104.1765
$ws.Range("K61").Value = 1104.1765
$ws.Range("M61").Value = -892.1765
$ws.Range("H74").Value = 612561.1
$ws.Range("I74").Value = 1501527.8
$ws.Range("J74").Value = 19916.666
$ws.Range("K74").Value = 1501527.8
$ws.Range("L74").Value = 19916.666
$ws.Range("M74").Value = -1500653.8
$ws.Range("N74").Value = -21664.666
$ws.Range("H77").Value = 612561.1
$ws.Range("I77").Value = 1501527.8
$ws.Range("J77").Value = 19916.666
$ws.Range("K77").Value = 7507639
$ws.Range("L77").Value = 99583.33
$ws.Range("M77").Value = -7503271
$ws.Range("N77").Value = -108319.33
$ws.Range("H102").Value = 3787.5715
$ws.Range("I102").Value = 4275.364
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 4275.364
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = -2653.364
$ws.Range("N102").Value = -5243
$ws.Range("H122").Value = 2587.3076
$ws.Range("I122").Value = 2536.25
$ws.Range("J122").Value = 2757.5
$ws.Range("K122").Value = 7608.75
$ws.Range("L122").Value = 8272.5
$ws.Range("M122").Value = -5158.75
$ws.Range("N122").Value = -13172.5
$ws.Range("H132").Value = 980.2273
$ws.Range("I132").Value = 795.5854
$ws.Range("K132").Value = 2386.7562
$ws.Range("M132").Value = 143.2437999999997
$ws.Range("H136").Value = 6303.136
$ws.Range("I136").Value = 1104.1765
$ws.Range("K136").Value = 3312.5295
$ws.Range("M136").Value = -762.5295000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 53587.668
$ws.Range("I20").Value = 53587.668
$ws.Range("K20").Value = 53587.668
$ws.Range("M20").Value = -53340.668
$ws.Range("H58").Value = 34835.25
$ws.Range("J58").Value = 42180.332
$ws.Range("L58").Value = 42180.332
$ws.Range("N58").Value = -42768.332
$ws.Range("H105").Value = 1901.2858
$ws.Range("I105").Value = 1461.8
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 1461.8
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = 285.2
$ws.Range("N105").Value = -6494

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3127635.2
$ws.Range("I31").Value = 5000998
$ws.Range("J31").Value = 5364.25
$ws.Range("K31").Value = 5000998
$ws.Range("L31").Value = 5364.25
$ws.Range("M31").Value = -5000703
$ws.Range("N31").Value = -5954.25
$ws.Range("H34").Value = 3127635.2
$ws.Range("I34").Value = 5000998
$ws.Range("J34").Value = 5364.25
$ws.Range("K34").Value = 5000998
$ws.Range("L34").Value = 5364.25
$ws.Range("M34").Value = -5000796
$ws.Range("N34").Value = -5768.25
$ws.Range("H58").Value = 17360.143
$ws.Range("I58").Value = 1959.1765
$ws.Range("K58").Value = 1959.1765
$ws.Range("M58").Value = -1756.1765
$ws.Range("H62").Value = 5483.8184
$ws.Range("I62").Value = 5710.2856
$ws.Range("J62").Value = 5087.5
$ws.Range("K62").Value = 5710.2856
$ws.Range("L62").Value = 5087.5
$ws.Range("M62").Value = -5086.2856
$ws.Range("N62").Value = -6335.5
$ws.Range("H65").Value = 5483.8184
$ws.Range("I65").Value = 5710.2856
$ws.Range("J65").Value = 5087.5
$ws.Range("K65").Value = 28551.428
$ws.Range("L65").Value = 25437.5
$ws.Range("M65").Value = -25431.428
$ws.Range("N65").Value = -31677.5
$ws.Range("H86").Value = 45790.42
$ws.Range("I86").Value = 59770.77
$ws.Range("J86").Value = 15499.667
$ws.Range("K86").Value = 59770.77
$ws.Range("L86").Value = 15499.667
$ws.Range("M86").Value = -58647.77
$ws.Range("N86").Value = -17745.667
$ws.Range("H89").Value = 45790.42
$ws.Range("I89").Value = 59770.77
$ws.Range("J89").Value = 15499.667
$ws.Range("K89").Value = 298853.85
$ws.Range("L89").Value = 77498.33499999999
$ws.Range("M89").Value = -293237.85
$ws.Range("N89").Value = -88730.33499999999
$ws.Range("H132").Value = 54064.367
$ws.Range("I132").Value = 63513.938
$ws.Range("K132").Value = 190541.814
$ws.Range("M132").Value = -188011.814
$ws.Range("H136").Value = 17360.143
$ws.Range("I136").Value = 1959.1765
$ws.Range("K136").Value = 5877.529500000001
$ws.Range("M136").Value = -3327.529500000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 5984
$ws.Range("I140").Value = 5984
$ws.Range("K140").Value = 17952
$ws.Range("M140").Value = -12772

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 196.73334
$ws.Range("I2").Value = 105.75
$ws.Range("K2").Value = 105.75
$ws.Range("M2").Value = 7.25
$ws.Range("H55").Value = 7600
$ws.Range("I55").Value = 4025
$ws.Range("J55").Value = 14750
$ws.Range("K55").Value = 4025
$ws.Range("L55").Value = 14750
$ws.Range("M55").Value = -3698
$ws.Range("N55").Value = -15404
$ws.Range("H80").Value = 25500.143
$ws.Range("I80").Value = 1175
$ws.Range("J80").Value = 35230.2
$ws.Range("K80").Value = 1175
$ws.Range("L80").Value = 35230.2
$ws.Range("M80").Value = -177
$ws.Range("N80").Value = -37226.2
$ws.Range("H83").Value = 25500.143
$ws.Range("I83").Value = 1175
$ws.Range("J83").Value = 35230.2
$ws.Range("K83").Value = 5875
$ws.Range("L83").Value = 176151
$ws.Range("M83").Value = -883
$ws.Range("N83").Value = -186135
$ws.Range("H92").Value = 18749
$ws.Range("J92").Value = 18749
$ws.Range("L92").Value = 18749
$ws.Range("N92").Value = -22493
$ws.Range("H122").Value = 4451.577
$ws.Range("J122").Value = 4793.5557
$ws.Range("L122").Value = 14380.6671
$ws.Range("N122").Value = -19280.6671
$ws.Range("H132").Value = 3776.4167
$ws.Range("I132").Value = 3670.375
$ws.Range("J132").Value = 4624.75
$ws.Range("K132").Value = 11011.125
$ws.Range("L132").Value = 13874.25
$ws.Range("M132").Value = -8481.125
$ws.Range("N132").Value = -18934.25
$ws.Range("H136").Value = 7666.2856
$ws.Range("J136").Value = 7666.2856
$ws.Range("L136").Value = 22998.8568
$ws.Range("N136").Value = -28098.8568

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 22500
$ws.Range("J4").Value = 22500
$ws.Range("L4").Value = 22500
$ws.Range("N4").Value = -22726
$ws.Range("H25").Value = 56669
$ws.Range("J25").Value = 80000
$ws.Range("L25").Value = 80000
$ws.Range("N25").Value = -80460
$ws.Range("H28").Value = 22500
$ws.Range("J28").Value = 22500
$ws.Range("L28").Value = 22500
$ws.Range("N28").Value = -22964
$ws.Range("H37").Value = 22500
$ws.Range("J37").Value = 22500
$ws.Range("L37").Value = 22500
$ws.Range("N37").Value = -22714

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 82688.336
$ws.Range("J75").Value = 82688.336
$ws.Range("L75").Value = 82688.336
$ws.Range("N75").Value = -84560.336
$ws.Range("H78").Value = 82688.336
$ws.Range("J78").Value = 82688.336
$ws.Range("L78").Value = 248065.008
$ws.Range("N78").Value = -257425.008
$ws.Range("H132").Value = 2236.818
$ws.Range("I132").Value = 1419.625
$ws.Range("K132").Value = 4258.875
$ws.Range("M132").Value = -1728.875
$ws.Range("H136").Value = 15732.195
$ws.Range("I136").Value = 16661.21
$ws.Range("J136").Value = 3964.6667
$ws.Range("K136").Value = 49983.63
$ws.Range("L136").Value = 11894.0001
$ws.Range("M136").Value = -47433.63
$ws.Range("N136").Value = -16994.0001
